# Weekly update: a new Brócoli price record (week of 2023-10-13) was
# reported for Terminal Hortofrutícola Agro Chillán. It is inserted as a
# new row 565, pushing all subsequent historical rows down by one
# (old row 565 -> 566, ..., old row 628 -> 629).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 565; everything below shifts down.
$ws.Rows("565:565").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(565, 1).Value = 7
$ws.Cells.Item(565, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(565, 3).Value = "Ñuble"
$ws.Cells.Item(565, 4).Value = 45212
$ws.Cells.Item(565, 5).Value = 16
$ws.Cells.Item(565, 6).Value = 100112023
$ws.Cells.Item(565, 7).Value = "Brócoli"
$ws.Cells.Item(565, 8).Value = "Sin especificar"
$ws.Cells.Item(565, 9).Value = "Primera"
$ws.Cells.Item(565, 10).Value = 500
$ws.Cells.Item(565, 11).Value = 1000
$ws.Cells.Item(565, 12).Value = 1000
$ws.Cells.Item(565, 13).Value = 1000
$ws.Cells.Item(565, 14).Value = "$/unidad"
$ws.Cells.Item(565, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(565, 16).Value = 1000
$ws.Cells.Item(565, 17).Value = 1
$ws.Cells.Item(565, 18).Value = "Hortaliza"
